$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# A new (blank) column is inserted before column N ("Late"), pushing the
# existing N/O/P columns ("Late", "heading", "Outstanding") one slot to the
# right (-> O/P/Q). Values are carried along automatically by the insert.
$ws.Columns("N").Insert()

# The newly inserted column picks up the same width as its neighbour
# (column M, "In Advance" = 10.7109375 chars) but without the bestFit flag,
# matching a manually-set (not autofit) width.
$ws.Columns("N").ColumnWidth = 9.877604166666666

# Make "Repayment schedule" the active/selected sheet (was "Transactions"),
# and move its selection to J15.
[void]$ws.Activate()
[void]$ws.Range("J15").Select()
